$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -230
$ws.Range("N12").ClearContents()
$ws.Range("H15").Value = 999.55316
$ws.Range("I15").Value = 999.55316
$ws.Range("K15").Value = 2998.65948
$ws.Range("M15").Value = -2829.65948
$ws.Range("H38").Value = 137.09525
$ws.Range("I38").Value = 88.36842
$ws.Range("J38").Value = 600
$ws.Range("K38").Value = 265.10526
$ws.Range("L38").Value = 1800
$ws.Range("M38").Value = 106.89474
$ws.Range("N38").Value = -2544
$ws.Range("H40").Value = 4896.241
$ws.Range("I40").Value = 8434.214
$ws.Range("J40").Value = 1594.1333
$ws.Range("K40").Value = 8434.214
$ws.Range("L40").Value = 1594.1333
$ws.Range("M40").Value = -8259.214
$ws.Range("N40").Value = -1944.1333
$ws.Range("H132").Value = 1292.6052
$ws.Range("I132").Value = 967.55884
$ws.Range("J132").Value = 4055.5
$ws.Range("K132").Value = 2902.67652
$ws.Range("L132").Value = 12166.5
$ws.Range("M132").Value = -372.67652
$ws.Range("N132").Value = -17226.5
$ws.Range("H137").Value = 1567.1333
$ws.Range("I137").Value = 1337
$ws.Range("K137").Value = 4011
$ws.Range("M137").Value = -1461

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 551.5333000000001
$ws.Range("I97").Value = 606.73914
$ws.Range("J97").Value = 370.14285
$ws.Range("K97").Value = 606.73914
$ws.Range("L97").Value = 370.14285
$ws.Range("M97").Value = -110.73914
$ws.Range("N97").Value = -1362.14285
$ws.Range("H132").Value = 2672.158
$ws.Range("I132").Value = 1432.8182
$ws.Range("J132").Value = 4376.25
$ws.Range("K132").Value = 4298.4546
$ws.Range("L132").Value = 13128.75
$ws.Range("M132").Value = -1768.4546
$ws.Range("N132").Value = -18188.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1787
$ws.Range("I31").Value = 870.069
$ws.Range("J31").Value = 5585.7144
$ws.Range("K31").Value = 870.069
$ws.Range("L31").Value = 5585.7144
$ws.Range("M31").Value = -575.069
$ws.Range("N31").Value = -6175.7144
$ws.Range("H34").Value = 1787
$ws.Range("I34").Value = 870.069
$ws.Range("J34").Value = 5585.7144
$ws.Range("K34").Value = 870.069
$ws.Range("L34").Value = 5585.7144
$ws.Range("M34").Value = -668.069
$ws.Range("N34").Value = -5989.7144
$ws.Range("H99").Value = 1800
$ws.Range("I99").Value = 1600
$ws.Range("K99").Value = 1600
$ws.Range("M99").Value = -102
$ws.Range("H105").Value = 1234.9375
$ws.Range("I105").Value = 896.2727
$ws.Range("K105").Value = 896.2727
$ws.Range("M105").Value = 850.7273
$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 1600
$ws.Range("K126").Value = 4800
$ws.Range("M126").Value = -2330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 397.375
$ws.Range("I5").Value = 397.375
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1192.125
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1080.125
$ws.Range("N5").ClearContents()
$ws.Range("H88").Value = 3911.4285
$ws.Range("J88").Value = 3911.4285
$ws.Range("L88").Value = 11734.2855
$ws.Range("N88").Value = -12590.2855
$ws.Range("H91").Value = 3911.4285
$ws.Range("J91").Value = 3911.4285
$ws.Range("L91").Value = 11734.2855
$ws.Range("N91").Value = -14698.2855
$ws.Range("H110").Value = 7055.7144
$ws.Range("I110").Value = 3433.3333
$ws.Range("J110").Value = 9772.5
$ws.Range("K110").Value = 10299.9999
$ws.Range("L110").Value = 29317.5
$ws.Range("M110").Value = -6209.999899999999
$ws.Range("N110").Value = -37497.5
$ws.Range("H122").Value = 6911.9375
$ws.Range("I122").Value = 465.77777
$ws.Range("J122").Value = 15199.857
$ws.Range("K122").Value = 4191.99993
$ws.Range("L122").Value = 136798.713
$ws.Range("M122").Value = -1741.99993
$ws.Range("N122").Value = -141698.713
$ws.Range("H135").Value = 397.375
$ws.Range("I135").Value = 397.375
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3576.375
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1041.375
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 289.23077
$ws.Range("I55").Value = 211.25
$ws.Range("J55").Value = 414
$ws.Range("K55").Value = 211.25
$ws.Range("L55").Value = 414
$ws.Range("M55").Value = -38.25
$ws.Range("N55").Value = -760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2425
$ws.Range("I62").Value = 2900
$ws.Range("J62").Value = 2266.6667
$ws.Range("K62").Value = 2900
$ws.Range("L62").Value = 2266.6667
$ws.Range("M62").Value = -2276
$ws.Range("N62").Value = -3514.6667
$ws.Range("H65").Value = 2425
$ws.Range("I65").Value = 2900
$ws.Range("J65").Value = 2266.6667
$ws.Range("K65").Value = 14500
$ws.Range("L65").Value = 11333.3335
$ws.Range("M65").Value = -11380
$ws.Range("N65").Value = -17573.3335
$ws.Range("H81").Value = 4774.0605
$ws.Range("I81").Value = 10941
$ws.Range("J81").Value = 2092.7827
$ws.Range("K81").Value = 21882
$ws.Range("L81").Value = 4185.5654
$ws.Range("M81").Value = -20821
$ws.Range("N81").Value = -6307.5654
$ws.Range("H84").Value = 4774.0605
$ws.Range("I84").Value = 10941
$ws.Range("J84").Value = 2092.7827
$ws.Range("K84").Value = 109410
$ws.Range("L84").Value = 20927.827
$ws.Range("M84").Value = -104106
$ws.Range("N84").Value = -31535.827
$ws.Range("H96").Value = 2316.077
$ws.Range("I96").Value = 2078.7778
$ws.Range("J96").Value = 2850
$ws.Range("K96").Value = 2078.7778
$ws.Range("L96").Value = 2850
$ws.Range("M96").Value = -705.7777999999998
$ws.Range("N96").Value = -5596
$ws.Range("H100").Value = 996.375
$ws.Range("I100").Value = 1103.2307
$ws.Range("K100").Value = 2206.4614
$ws.Range("M100").Value = -1665.4614
$ws.Range("H107").Value = 4445.778
$ws.Range("I107").Value = 898.93335
$ws.Range("J107").Value = 8879.333000000001
$ws.Range("K107").Value = 2696.80005
$ws.Range("L107").Value = 26637.999
$ws.Range("M107").Value = -776.8000499999998
$ws.Range("N107").Value = -30477.999
$ws.Range("H113").Value = 259
$ws.Range("I113").Value = 160.375
$ws.Range("J113").Value = 357.625
$ws.Range("K113").Value = 481.125
$ws.Range("L113").Value = 1072.875
$ws.Range("M113").Value = 1688.875
$ws.Range("N113").Value = -5412.875
$ws.Range("H122").Value = 1479.7858
$ws.Range("I122").Value = 1161.7727
$ws.Range("J122").Value = 2645.8333
$ws.Range("K122").Value = 3485.3181
$ws.Range("L122").Value = 7937.499899999999
$ws.Range("M122").Value = -1035.3181
$ws.Range("N122").Value = -12837.4999
$ws.Range("H126").Value = 1033.875
$ws.Range("I126").Value = 839.8182
$ws.Range("J126").Value = 1460.8
$ws.Range("K126").Value = 2519.4546
$ws.Range("L126").Value = 4382.4
$ws.Range("M126").Value = -49.45460000000003
$ws.Range("N126").Value = -9322.4
$ws.Range("H136").Value = 4048.6572
$ws.Range("I136").Value = 881.6667
$ws.Range("J136").Value = 7401.9414
$ws.Range("K136").Value = 2645.0001
$ws.Range("L136").Value = 22205.8242
$ws.Range("M136").Value = -95.0001000000002
$ws.Range("N136").Value = -27305.8242
